# fall 22 week 12 complete
# Appends 24 new rows of matchup data (rows 1433-1456) to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
    @(4, 2, 4, 0),
    @(4, 2, 3, 1),
    @(6, 0, 7, 3),
    @(4, 2, 5, 0),
    @(5, 2, 3, 1),
    @(3, 2, 3, 1),
    @(6, 1, 6, 2),
    @(3, 1, 5, 2),
    @(7, 3, 6, 0),
    @(3, 3, 3, 0),
    @(5, 0, 4, 2),
    @(6, 2, 3, 1),
    @(3, 3, 3, 0),
    @(3, 0, 3, 3),
    @(3, 3, 3, 0),
    @(5, 2, 4, 1),
    @(4, 2, 4, 0),
    @(5, 2, 5, 1),
    @(2, 0, 6, 3),
    @(5, 0, 3, 2),
    @(5, 2, 5, 0),
    @(5, 2, 6, 0),
    @(4, 1, 4, 2),
    @(4, 0, 3, 2)
)

$startRow = 1433
$endRow = $startRow + $data.Count - 1

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# Update view to reflect the new active cell selection (mirrors Excel
# auto-advancing the selection past the last entered row).
$nextRow = $endRow + 1
$ws.Activate()
$ws.Range("A$nextRow").Select()
